$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Price" column (D) - force text format to preserve exact
# formatting (trailing zeros, multi-dot thousand separators, etc.)
# that Excel would otherwise mangle by auto-converting to a number.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.926.10'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.890.20'
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7254'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.96'
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3087'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '26.04'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06874'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07943'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7660'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.885.97'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.232'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '90.90'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.996.69'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.07'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.730'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007738'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '238.31'
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.135.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.004'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.821'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.259'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.90'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.90'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1271'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.007'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.357'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.531'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.285'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.053'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05047'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.266'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7305'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.732'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01913'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.769'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.307'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.20'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4418'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.914'
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8349'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '100.57'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.563'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.737'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '37.38'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.039.76'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '938.50'
$ws.Range("D51").Style = "Normal"

# Update "Volume(1h)" column (E) - plain text assignment (the
# padding spaces already keep Excel from treating these as numbers).
$ws.Range("E2").Value = '  -1.55%  '
$ws.Range("E3").Value = '  -2.55%  '
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("E5").Value = '  -6.60%  '
$ws.Range("E6").Value = '  -2.03%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("E8").Value = '  -3.88%  '
$ws.Range("E9").Value = '  -6.59%  '
$ws.Range("E11").Value = '  -1.07%  '
$ws.Range("E12").Value = '  -2.20%  '
$ws.Range("E13").Value = '  -2.48%  '
$ws.Range("E14").Value = '  -2.67%  '
$ws.Range("E15").Value = '  -4.33%  '
$ws.Range("E16").Value = '  -1.28%  '
$ws.Range("E17").Value = '  -3.05%  '
$ws.Range("E18").Value = '  -1.55%  '
$ws.Range("E19").Value = '  -2.99%  '
$ws.Range("E20").Value = '  -6.94%  '
$ws.Range("E21").Value = '  +0.14%  '
$ws.Range("E22").Value = '  -2.43%  '
$ws.Range("E23").Value = '  +0.35%  '
$ws.Range("E24").Value = '  +0.86%  '
$ws.Range("E25").Value = '  -3.30%  '
$ws.Range("E26").Value = '  +0.72%  '
$ws.Range("E27").Value = '  -1.02%  '
$ws.Range("E28").Value = '  -6.30%  '
$ws.Range("E29").Value = '  -12.09%  '
$ws.Range("E30").Value = '  -1.02%  '
$ws.Range("E31").Value = '  +0.60%  '
$ws.Range("E32").Value = '  -3.29%  '
$ws.Range("E33").Value = '  -1.85%  '
$ws.Range("E34").Value = '  -2.41%  '
$ws.Range("E35").Value = '  -1.49%  '
$ws.Range("E36").Value = '  -2.45%  '
$ws.Range("E37").Value = '  -1.82%  '
$ws.Range("E38").Value = '  -2.46%  '
$ws.Range("E39").Value = '  -1.37%  '
$ws.Range("E40").Value = '  -2.28%  '
$ws.Range("E41").Value = '  -5.97%  '
$ws.Range("E42").Value = '  -2.32%  '
$ws.Range("E43").Value = '  -3.08%  '
$ws.Range("E44").Value = '  -0.34%  '
$ws.Range("E45").Value = '  -0.38%  '
$ws.Range("E46").Value = '  -0.38%  '
$ws.Range("E47").Value = '  +0.50%  '
$ws.Range("E48").Value = '  -0.05%  '
$ws.Range("E49").Value = '  +0.33%  '
$ws.Range("E50").Value = '  -2.21%  '
$ws.Range("E51").Value = '  -5.18%  '
